$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("A7:I7")
$src.Copy($ws.Range("A8:I8"))
$src.Copy($ws.Range("A9:I9"))
$src.Copy($ws.Range("A10:I10"))
